# Update monthly submission report ("testing.xlsx") for the new month's
# submissions: bump January 2017 (column I) counts on the "Sheet" tab,
# fill in a few previously-blank Name/Email pairs, append 10 new
# submitter rows (227-236), and extend the "Monthly_STAT" SUM() formulas
# so they cover the new row range (2:236 instead of 2:226).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$stat = $wb.Worksheets.Item("Monthly_STAT")

function Set-Cell {
    param($sheet, $row, $col, $value)
    $c = $sheet.Cells.Item($row, $col)
    $c.Value = $value
    # Keep formatting identical to the untouched cells around it (no
    # explicit per-cell style override) instead of inheriting the
    # column's default style id.
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# Column letters -> numbers for the "Sheet" tab: A=1 UID, B=2 Name,
# C=3 Email, D=4 Aug2016 .. I=9 Jan2017
# ---------------------------------------------------------------------

# -- Existing-row updates (mostly January 2017 / column I counts) -----
Set-Cell $ws 6   9 2
Set-Cell $ws 11  9 2
Set-Cell $ws 12  9 1
Set-Cell $ws 14  9 115
Set-Cell $ws 15  9 1
Set-Cell $ws 21  9 4
Set-Cell $ws 23  9 17
Set-Cell $ws 32  9 1
Set-Cell $ws 37  9 1
Set-Cell $ws 40  9 2
Set-Cell $ws 70  9 1
Set-Cell $ws 72  9 1
Set-Cell $ws 79  9 1
Set-Cell $ws 96  9 1
Set-Cell $ws 100 9 1
Set-Cell $ws 102 9 1
Set-Cell $ws 110 9 3
Set-Cell $ws 113 9 2
Set-Cell $ws 141 9 12
Set-Cell $ws 155 9 4
Set-Cell $ws 158 9 12
Set-Cell $ws 175 9 1
Set-Cell $ws 182 9 5

# Row 183 (UID CFD8AFA4C0) gains a resolved Name/Email plus a count
Set-Cell $ws 183 2 "arango juan"
Set-Cell $ws 183 3 "arango.juan@jjay.cuny.edu"
Set-Cell $ws 183 9 3

Set-Cell $ws 184 9 1
Set-Cell $ws 189 9 2
Set-Cell $ws 194 9 3
Set-Cell $ws 197 9 1
Set-Cell $ws 199 9 2

# -- New submitter rows (227-236) --------------------------------------
Set-Cell $ws 227 1 "0FD8B42CC0"
Set-Cell $ws 227 2 "amin shahid"
Set-Cell $ws 227 3 "amin.shahid@jjay.cuny.edu"
Set-Cell $ws 227 9 3

Set-Cell $ws 228 1 "8FD8A3A820"
Set-Cell $ws 228 2 "cowell,truman"
Set-Cell $ws 228 3 "cowell.truman@jjay.cuny.edu"
Set-Cell $ws 228 9 3

Set-Cell $ws 229 1 "8FD8B367A0"
Set-Cell $ws 229 9 1

Set-Cell $ws 230 1 "4FD8B36A40"
Set-Cell $ws 230 2 "Jean Augustin"
Set-Cell $ws 230 3 "Jean.Augustine@jjay.cuny.edu"
Set-Cell $ws 230 9 2

Set-Cell $ws 231 1 "0FD8B45C20"
Set-Cell $ws 231 9 1

Set-Cell $ws 232 1 "8FD8B68DE0"
Set-Cell $ws 232 2 "Miguel Martillo"
Set-Cell $ws 232 3 "miguel.martillo@jjay.cuny.edu"
Set-Cell $ws 232 9 3

Set-Cell $ws 233 1 "8FD8B65040"
Set-Cell $ws 233 9 1

Set-Cell $ws 234 1 "CFD8AEC840"
Set-Cell $ws 234 9 1

Set-Cell $ws 235 1 "8FD8ADBD20"
Set-Cell $ws 235 9 1

Set-Cell $ws 236 1 "8FD8A91340"
Set-Cell $ws 236 2 "Jesse Silkworth"
Set-Cell $ws 236 3 "jesse,silkworth@jjay.cuny.edu"
Set-Cell $ws 236 9 1

# ---------------------------------------------------------------------
# Monthly_STAT: extend each SUM() range from row 226 to row 236
# ---------------------------------------------------------------------
$stat.Range("B2").Formula = "=SUM(Sheet!D2:D236)"
$stat.Range("B3").Formula = "=SUM(Sheet!E2:E236)"
$stat.Range("B4").Formula = "=SUM(Sheet!F2:F236)"
$stat.Range("B5").Formula = "=SUM(Sheet!G2:G236)"
$stat.Range("B6").Formula = "=SUM(Sheet!H2:H236)"
$stat.Range("B7").Formula = "=SUM(Sheet!I2:I236)"
